{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the \"correction processus de calcul facturation\" changes:\n//   - Both \"Date de cr\u00e9ation\" / \"Date de version\" fields: 15/06/2018 -> 18/06/2018\n//   - Line item amount 6000,00\u20ac -> 6220,00\u20ac\n//   - Line item amount 4890,00\u20ac -> 4975,00\u20ac\n//   - Sub-total 25155,00 -> 25460,00\n//   - Grand total 61420,00 -> 61725,00\n\nconst replacements = [\n  { from: \"15/06/2018\", to: \"18/06/2018\" },\n  { from: \"15/06/2018\", to: \"18/06/2018\" },\n  { from: \"6000,00\u20ac\", to: \"6220,00\u20ac\" },\n  { from: \"4890,00\u20ac\", to: \"4975,00\u20ac\" },\n  { from: \"25155,00\", to: \"25460,00\" },\n  { from: \"61420,00\", to: \"61725,00\" },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${from}`);\n  }\n\n  // Only the first (not-yet-updated) occurrence should be changed per call,\n  // so repeated \"from\" values (e.g. the two identical dates) are each\n  // consumed exactly once across iterations.\n  results.items[0].insertText(to, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the \"correction processus de calcul facturation\" changes:\n#   - Both \"Date de cr\u00e9ation\" / \"Date de version\" fields: 15/06/2018 -> 18/06/2018\n#   - Line item amount 6000,00\u20ac -> 6220,00\u20ac\n#   - Line item amount 4890,00\u20ac -> 4975,00\u20ac\n#   - Sub-total 25155,00 -> 25460,00\n#   - Grand total 61420,00 -> 61725,00\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n\nfunction Replace-Once($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, $wdReplaceOne)\n}\n\nReplace-Once \"15/06/2018\" \"18/06/2018\"\nReplace-Once \"15/06/2018\" \"18/06/2018\"\nReplace-Once \"6000,00\u20ac\" \"6220,00\u20ac\"\nReplace-Once \"4890,00\u20ac\" \"4975,00\u20ac\"\nReplace-Once \"25155,00\" \"25460,00\"\nReplace-Once \"61420,00\" \"61725,00\"\n"}
